$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

function Set-Row {
    param($ws, $r, $a, $b, $c, $d, $e)
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# --- Update header metadata (Última actualización / Total filas) for all 3 sheets ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 14:00:52"
$ws1.Cells.Item(3,1).Value = "Total filas: 268"

$ws2.Cells.Item(2,1).Value = "Última actualización: 14:00:52"
$ws2.Cells.Item(3,1).Value = "Total filas: 28"

$ws3.Cells.Item(2,1).Value = "Última actualización: 14:00:52"
$ws3.Cells.Item(3,1).Value = "Total filas: 37"

# --- Sheet1 (LP1912): tie-break reorderings earlier in the table ---
Set-Row $ws1 71 "06:43:12" "08:23" "16_P MOR-SANTA ANA" 100 "LP1912"
Set-Row $ws1 72 "06:43:12" "08:23" "215B_EL PATO" 100 "LP1912"

Set-Row $ws1 96 "08:36:20" "09:11" "16_SANTA ANA" 35 "LP1912"
Set-Row $ws1 97 "08:48:09" "09:11" "16_P MOR-SANTA ANA" 23 "LP1912"

Set-Row $ws1 105 "07:58:19" "09:23" "17_ROMERO" 85 "LP1912"
Set-Row $ws1 106 "07:45:49" "09:23" "11_ETCHEVERRY" 98 "LP1912"

# --- Sheet1 (LP1912): new scraped rows appended/merged from row 240 onward ---
Set-Row $ws1 240 "14:00:52" "14:00" "16_SANTA ANA" 0 "LP1912"
Set-Row $ws1 241 "14:00:52" "14:00" "14_ABASTO" 0 "LP1912"
Set-Row $ws1 242 "14:00:52" "14:01" "15_ABASTO" 1 "LP1912"
Set-Row $ws1 243 "14:00:52" "14:04" "23_HERNANDEZ" 4 "LP1912"
Set-Row $ws1 244 "12:27:08" "14:04" "17_ROMERO" 97 "LP1912"
Set-Row $ws1 245 "13:23:09" "14:05" "23_HERNANDEZ" 42 "LP1912"
Set-Row $ws1 246 "14:00:52" "14:06" "16_SANTA ANA" 6 "LP1912"
Set-Row $ws1 247 "14:00:52" "14:16" "27_EL RETIRO" 16 "LP1912"
Set-Row $ws1 248 "12:27:08" "14:17" "27_EL RETIRO" 110 "LP1912"
Set-Row $ws1 249 "14:00:52" "14:19" "215C_EL PATO" 19 "LP1912"
Set-Row $ws1 250 "14:00:52" "14:20" "26_HERNANDEZ" 20 "LP1912"
Set-Row $ws1 251 "12:27:08" "14:20" "215C_EL PATO" 113 "LP1912"
Set-Row $ws1 252 "12:54:06" "14:21" "26_HERNANDEZ" 87 "LP1912"
Set-Row $ws1 253 "12:54:06" "14:39" "14_ABASTO" 105 "LP1912"
Set-Row $ws1 254 "14:00:52" "14:44" "14_ABASTO" 44 "LP1912"
Set-Row $ws1 255 "14:00:52" "14:56" "16_P MOR-SANTA ANA" 56 "LP1912"
Set-Row $ws1 256 "13:23:09" "14:57" "16_P MOR-SANTA ANA" 94 "LP1912"
Set-Row $ws1 257 "13:23:09" "14:58" "215B_EL PATO" 95 "LP1912"
Set-Row $ws1 258 "13:23:09" "15:00" "81_EL PELIGRO" 97 "LP1912"
Set-Row $ws1 259 "14:00:52" "15:04" "10_OLMOS" 64 "LP1912"
Set-Row $ws1 260 "13:23:09" "15:05" "10_OLMOS" 102 "LP1912"
Set-Row $ws1 261 "14:00:52" "15:10" "17_ROMERO" 70 "LP1912"
Set-Row $ws1 262 "14:00:52" "15:13" "11_ETCHEVERRY" 73 "LP1912"
Set-Row $ws1 263 "13:23:09" "15:14" "11_ETCHEVERRY" 111 "LP1912"
Set-Row $ws1 264 "14:00:52" "15:20" "15_ABASTO" 80 "LP1912"
Set-Row $ws1 265 "13:23:09" "15:21" "26_HERNANDEZ" 118 "LP1912"
Set-Row $ws1 266 "14:00:52" "15:25" "26_HERNANDEZ" 85 "LP1912"
Set-Row $ws1 267 "14:00:52" "15:32" "84_COLONIA URQUIZA-ESC 49" 92 "LP1912"
Set-Row $ws1 268 "14:00:52" "15:35" "23_HERNANDEZ" 95 "LP1912"
Set-Row $ws1 269 "14:00:52" "15:36" "10_OLMOS" 96 "LP1912"
Set-Row $ws1 270 "14:00:52" "15:38" "215A_EL PATO" 98 "LP1912"
Set-Row $ws1 271 "14:00:52" "15:46" "16_P MOR-167 Y 521" 106 "LP1912"
Set-Row $ws1 272 "14:00:52" "15:53" "11_ETCHEVERRY" 113 "LP1912"
Set-Row $ws1 273 "14:00:52" "15:56" "27_EL RETIRO" 116 "LP1912"

# --- Sheet2 (LP1912-215): new scraped rows appended/merged from row 30 onward ---
Set-Row $ws2 30 "14:00:52" "14:19" "215C_EL PATO" 19 "LP1912"
Set-Row $ws2 31 "12:27:08" "14:20" "215C_EL PATO" 113 "LP1912"
Set-Row $ws2 32 "13:23:09" "14:58" "215B_EL PATO" 95 "LP1912"
Set-Row $ws2 33 "14:00:52" "15:38" "215A_EL PATO" 98 "LP1912"

# --- Sheet3 (6203-6173): new scraped rows appended/merged from row 40 onward ---
Set-Row $ws3 40 "14:00:52" "14:52" "215D_LA PLATA" 52 "L6203"
Set-Row $ws3 41 "12:54:06" "14:53" "215D_LA PLATA" 119 "L6203"
Set-Row $ws3 42 "14:00:52" "15:34" "215A_LA PLATA" 94 "L6173"
